# "Generate Report for Archive"
#
# The localization-status report is regenerated: the status text for the
# 166d75bf-...md item moves from "Ready for handoff" to "In Translation"
# on every sheet that shows it (Overview + the per-locale zh-cn / de-de
# sheets), and the now-narrower status column is resized to fit the new
# (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) show the status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# --- Per-locale sheets: column C is the "Status" column ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Resize the status columns to fit the new, shorter text ---
# (ColumnWidth is expressed in characters of the Normal style font and is
# quantized by Excel's column-width grid; 12.5 lands on the closest grid
# step to the refreshed report's column width.)
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
